$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text updates (Volume number + report date range) ----
$ws.Range("A8").Value = "Volume 30   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/16/2023  Through  1/22/2023"

# ---- Pass 1: force numeric-looking target strings ('0') to stay text ----
# (set format to Text, assign the value, then the later format-repaste in Pass 2
#  restores the real target style; order matters: this must run first)
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"

# ---- Pass 2: copy formats from donor cells (one cell at a time; multi-area PasteSpecial is unreliable) ----
foreach ($ref in @("M14","L15","N15","M22","E26","L26","M28","M29")) {
    $ws.Range("K22").Copy()
    $ws.Range($ref).PasteSpecial(-4122)
}
foreach ($ref in @("C15","C26","C27","D27","E27")) {
    $ws.Range("A14").Copy()
    $ws.Range($ref).PasteSpecial(-4122)
}
foreach ($ref in @("C16","D26","C30","F30","I30")) {
    $ws.Range("F15").Copy()
    $ws.Range($ref).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# ---- Pass 3: set final values for style-changed cells (skip the ones already set as protected text in Pass 1) ----
$ws.Range("M14").Value = -100
$ws.Range("L15").Value = -50
$ws.Range("N15").Value = 0
$ws.Range("M22").Value = 200
$ws.Range("E26").Value = -100
$ws.Range("L26").Value = -50
$ws.Range("M28").Value = -100
$ws.Range("M29").Value = -100
$ws.Range("E27").Value = "***.*"
$ws.Range("C16").Value = 2
$ws.Range("D26").Value = 2
$ws.Range("C30").Value = 1
$ws.Range("F30").Value = 1
$ws.Range("I30").Value = 1

# ---- Value-only cells (style already correct) ----
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = -63.157894736842
$ws.Range("I16").Value = 6
$ws.Range("J16").Value = 14
$ws.Range("K16").Value = -57.142857142857
$ws.Range("M16").Value = -79.310344827586
$ws.Range("N16").Value = -93.877551020408
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 60
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 142.857142857143
$ws.Range("I17").Value = 26
$ws.Range("J17").Value = 9
$ws.Range("K17").Value = 188.888888888889
$ws.Range("L17").Value = 52.941176470588
$ws.Range("M17").Value = 136.363636363636
$ws.Range("N17").Value = 30
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = -17.647058823529
$ws.Range("I18").Value = 11
$ws.Range("J18").Value = 9
$ws.Range("K18").Value = 22.222222222222
$ws.Range("L18").Value = 57.142857142857
$ws.Range("M18").Value = -31.25
$ws.Range("N18").Value = -89.523809523809
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 50
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = -27.906976744186
$ws.Range("I19").Value = 22
$ws.Range("J19").Value = 30
$ws.Range("K19").Value = -26.666666666666
$ws.Range("L19").Value = 69.230769230769
$ws.Range("N19").Value = -43.589743589743
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -57.142857142857
$ws.Range("F20").Value = 28
$ws.Range("G20").Value = 26
$ws.Range("H20").Value = 7.692307692307
$ws.Range("I20").Value = 19
$ws.Range("J20").Value = 21
$ws.Range("K20").Value = -9.523809523809
$ws.Range("L20").Value = 171.428571428571
$ws.Range("M20").Value = 35.714285714285
$ws.Range("N20").Value = -93.189964157706
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -16.666666666666
$ws.Range("F21").Value = 115
$ws.Range("G21").Value = 120
$ws.Range("H21").Value = -4.166666666666
$ws.Range("I21").Value = 85
$ws.Range("J21").Value = 84
$ws.Range("K21").Value = 1.190476190476
$ws.Range("L21").Value = 70
$ws.Range("M21").Value = 3.658536585365
$ws.Range("N21").Value = -84.317343173431
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 2
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -25
$ws.Range("I22").Value = 3
$ws.Range("J22").Value = 3
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 29.166666666666
$ws.Range("F24").Value = 114
$ws.Range("G24").Value = 104
$ws.Range("H24").Value = 9.615384615384
$ws.Range("I24").Value = 76
$ws.Range("J24").Value = 70
$ws.Range("K24").Value = 8.571428571428
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 58.333333333333
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 39
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = 34.482758620689
$ws.Range("I25").Value = 28
$ws.Range("J25").Value = 21
$ws.Range("K25").Value = 33.333333333333
$ws.Range("L25").Value = 47.368421052631
$ws.Range("M25").Value = -37.777777777777
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = -66.666666666666
$ws.Range("J26").Value = 3
$ws.Range("K26").Value = -66.666666666666
